{"js": "// Fix flaky test in user_resume_doc_spec.rb\n//\n// 1. The \"job description\" placeholder paragraph (a lone run whose text is\n//    exactly \"text\", sitting right after the Companies/Projects job entry)\n//    becomes \"_job_description_\".\n// 2. The stray placeholder paragraph whose text is exactly \"text.\" (right\n//    after the Conferences Talks entry, just before \"OpenSource\n//    Contributions\") is removed entirely.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet renamed = false;\nlet removedParagraph = null;\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n  if (!renamed && text === \"text\") {\n    paragraph.insertText(\"_job_description_\", \"Replace\");\n    renamed = true;\n  } else if (!removedParagraph && text === \"text.\") {\n    removedParagraph = paragraph;\n  }\n}\n\nif (removedParagraph) {\n  removedParagraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# Fix flaky test in user_resume_doc_spec.rb\n#\n# 1. The \"job description\" placeholder paragraph (a lone run whose text is\n#    exactly \"text\", right after the Companies/Projects job entry) becomes\n#    \"_job_description_\".\n# 2. The stray placeholder paragraph whose text is exactly \"text.\" (right\n#    after the Conferences Talks entry, just before \"OpenSource\n#    Contributions\") is removed entirely.\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13)\n    if ($t -eq \"text\") {\n        $p.Range.Text = \"_job_description_\"\n    } elseif ($t -eq \"text.\") {\n        $p.Range.Delete()\n    }\n}\n"}
